{"js": "// The \"Superficie\" (surface area) figure in the SHD description paragraph\n// changes from \"60,30\" to \"57,40\" (m2.). Everything else in the sentence\n// (\"... con una Superficie de 60,30 m2. en el inmueble ...\") stays the same.\nconst body = context.document.body;\n\nconst results = body.search(\"60,30\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Expected to find \"60,30\" in the document body, but it was not found.');\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"57,40\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The \"Superficie\" (surface area) figure in the SHD description paragraph\n# changes from \"60,30\" to \"57,40\" (m2.). Everything else in the sentence\n# (\"... con una Superficie de 60,30 m2. en el inmueble ...\") stays the same.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"60,30\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"57,40\"\n\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    throw 'Expected to find \"60,30\" in the document, but it was not found.'\n}\n"}
